# Y4_B2526_General_Surgery_checklist — "Log Time" column correction
# The D2:D74 "Log Time" values are corrected from ~11:53:34 AM to the
# accurate 11:03:15 AM, re-formatted with an explicit black-font h:mm:ss
# style, and the sheet selection is moved onto the corrected column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log-time value: 11:03:15 (39795 seconds past midnight).
$newTime = 39795.0 / 86400.0

# Set up the target style on the first data cell (D2): a plain h:mm:ss
# time format with an explicit black font, then propagate that exact
# style (and the corrected value) down through D74 without re-deriving
# it cell by cell, so only a single new style gets minted.
$first = $ws.Range("D2")
$first.NumberFormat = "h:mm:ss"
$first.Font.Color = 0
$first.Value = $newTime

$first.Copy()
$rest = $ws.Range("D3:D74")
$rest.PasteSpecial(-4122)  # xlPasteFormats
$rest.Value = $newTime

$excel.CutCopyMode = 0

# Reflect the edited column in the sheet's active selection.
$ws.Range("D2:D74").Select()
